$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Fix the AddCustomerTest row that previously had Runmode=N -> Y
$ws2.Range("A4").Value = "Y"

# Add the new "browser" column header and values for the existing 3 rows
$ws2.Range("E2").Value = "browser"
$ws2.Range("E3").Value = "chrome"
$ws2.Range("E4").Value = "firefox"

# Insert 8 new rows (5-12) before the old row 6 ("OpenAccountTest" block),
# pushing the OpenAccountTest section down from rows 6-12 to rows 14-20.
$ws2.Range("A5:A12").EntireRow.Insert()

# Fill the newly inserted rows with the repeated AddCustomerTest data,
# executing the test case across 10 different nodes (rows 3-12).
$ws2.Range("A5").Value = "Y"
$ws2.Range("B5").Value = "Ivan"
$ws2.Range("C5").Value = "Ivanov"
$ws2.Range("D5").Value = "e3r4t5"
$ws2.Range("E5").Value = "chrome"

$ws2.Range("A6").Value = "Y"
$ws2.Range("B6").Value = "Petr"
$ws2.Range("C6").Value = "Petrov"
$ws2.Range("D6").Value = "2af4g5"
$ws2.Range("E6").Value = "firefox"

$ws2.Range("A7").Value = "Y"
$ws2.Range("B7").Value = "Ivan"
$ws2.Range("C7").Value = "Ivanov"
$ws2.Range("D7").Value = "e3r4t5"
$ws2.Range("E7").Value = "chrome"

$ws2.Range("A8").Value = "Y"
$ws2.Range("B8").Value = "Petr"
$ws2.Range("C8").Value = "Petrov"
$ws2.Range("D8").Value = "2af4g5"
$ws2.Range("E8").Value = "firefox"

$ws2.Range("A9").Value = "Y"
$ws2.Range("B9").Value = "Ivan"
$ws2.Range("C9").Value = "Ivanov"
$ws2.Range("D9").Value = "e3r4t5"
$ws2.Range("E9").Value = "chrome"

$ws2.Range("A10").Value = "Y"
$ws2.Range("B10").Value = "Petr"
$ws2.Range("C10").Value = "Petrov"
$ws2.Range("D10").Value = "2af4g5"
$ws2.Range("E10").Value = "firefox"

$ws2.Range("A11").Value = "Y"
$ws2.Range("B11").Value = "Ivan"
$ws2.Range("C11").Value = "Ivanov"
$ws2.Range("D11").Value = "e3r4t5"
$ws2.Range("E11").Value = "chrome"

$ws2.Range("A12").Value = "Y"
$ws2.Range("B12").Value = "Petr"
$ws2.Range("C12").Value = "Petrov"
$ws2.Range("D12").Value = "2af4g5"
$ws2.Range("E12").Value = "firefox"

# Make TestData the active sheet/tab, and select the newly added block.
$ws2.Activate()
$ws2.Range("A11:E12").Select()
